# Populate the "07-05-22" sheet with the employee table (columns A:F, rows 1-3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("07-05-22")

# Header row
$ws.Range("A1").Value = "Фамилия"
$ws.Range("B1").Value = "Имя"
$ws.Range("C1").Value = "Отчество"
$ws.Range("D1").Value = "Должность"
$ws.Range("E1").Value = "Приоритет"
$ws.Range("F1").Value = "Пароль"

# Row 2
$ws.Range("A2").Value = "Полинкин"
$ws.Range("B2").Value = "Алексей"
$ws.Range("C2").Value = "Павлович"
$ws.Range("D2").Value = "Технический писатель"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "AlekseyPolinkin"

# Row 3
$ws.Range("A3").Value = "Кореньков"
$ws.Range("B3").Value = "Алексей"
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = "Конструктор"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "AlekseyKorenkov"

# Clear old leftover rows 4 and 5 (previous data had 5 rows, new table only has 3)
$ws.Range("A4:F5").ClearContents()

# Column widths to match target layout
# (target stored widths are 36.42578125 / 14 / 13.5703125 / 24.140625 / 14.140625 / 16.85546875;
#  ColumnWidth values below are chosen so the saved "width" attribute lands as close as possible)
$ws.Columns.Item(1).ColumnWidth = 35.666666666666664
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 23.333333333333332
$ws.Columns.Item(5).ColumnWidth = 13.333333333333334
$ws.Columns.Item(6).ColumnWidth = 16

# Update the selection anchor to match the saved view state
$ws.Range("F10").Select()
